# Update the regression-table "Estimate" / "Std. Error" figures in the
# single results table. Each replacement is scoped to Table.Cell(row,col)
# and uses wdReplaceOne (the 11th Find.Execute arg = 1) rather than
# wdReplaceAll, since a couple of values ("1.216") repeat verbatim in
# adjacent cells and a document-wide ReplaceAll would clobber both.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(2, 2)
$cell.Range.Find.Execute("156.908", $true, $false, $false, $false, $false, $true, 1, $false, "156.913", 1) | Out-Null

$cell = $t.Cell(2, 3)
$cell.Range.Find.Execute("0.475", $true, $false, $false, $false, $false, $true, 1, $false, "0.416", 1) | Out-Null

$cell = $t.Cell(3, 2)
$cell.Range.Find.Execute("117.430", $true, $false, $false, $false, $false, $true, 1, $false, "101.841", 1) | Out-Null

$cell = $t.Cell(3, 3)
$cell.Range.Find.Execute("1.216", $true, $false, $false, $false, $false, $true, 1, $false, "1.230", 1) | Out-Null

$cell = $t.Cell(4, 2)
$cell.Range.Find.Execute("114.015", $true, $false, $false, $false, $false, $true, 1, $false, "98.460", 1) | Out-Null

$cell = $t.Cell(4, 3)
$cell.Range.Find.Execute("1.216", $true, $false, $false, $false, $false, $true, 1, $false, "1.229", 1) | Out-Null

$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("103.781", $true, $false, $false, $false, $false, $true, 1, $false, "88.864", 1) | Out-Null

$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("1.336", $true, $false, $false, $false, $false, $true, 1, $false, "1.356", 1) | Out-Null

$cell = $t.Cell(6, 2)
$cell.Range.Find.Execute("140.034", $true, $false, $false, $false, $false, $true, 1, $false, "141.041", 1) | Out-Null

$cell = $t.Cell(6, 3)
$cell.Range.Find.Execute("1.047", $true, $false, $false, $false, $false, $true, 1, $false, "0.930", 1) | Out-Null

$cell = $t.Cell(7, 2)
$cell.Range.Find.Execute("121.194", $true, $false, $false, $false, $false, $true, 1, $false, "105.908", 1) | Out-Null

$cell = $t.Cell(7, 3)
$cell.Range.Find.Execute("2.618", $true, $false, $false, $false, $false, $true, 1, $false, "2.652", 1) | Out-Null

$cell = $t.Cell(8, 2)
$cell.Range.Find.Execute("116.892", $true, $false, $false, $false, $false, $true, 1, $false, "101.794", 1) | Out-Null

$cell = $t.Cell(8, 3)
$cell.Range.Find.Execute("2.651", $true, $false, $false, $false, $false, $true, 1, $false, "2.687", 1) | Out-Null

$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("106.658", $true, $false, $false, $false, $false, $true, 1, $false, "92.197", 1) | Out-Null

$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("2.675", $true, $false, $false, $false, $false, $true, 1, $false, "2.712", 1) | Out-Null
